$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3406
$ws.Range("I62").Value = 2933
$ws.Range("K62").Value = 2933
$ws.Range("M62").Value = -2309
$ws.Range("H64").Value = 5036.303
$ws.Range("I64").Value = 3528.842
$ws.Range("J64").Value = 7082.143
$ws.Range("K64").Value = 3528.842
$ws.Range("L64").Value = 7082.143
$ws.Range("M64").Value = -3280.842
$ws.Range("N64").Value = -7578.143
$ws.Range("H65").Value = 3406
$ws.Range("I65").Value = 2933
$ws.Range("K65").Value = 14665
$ws.Range("M65").Value = -11545
$ws.Range("H67").Value = 5036.303
$ws.Range("I67").Value = 3528.842
$ws.Range("J67").Value = 7082.143
$ws.Range("K67").Value = 3528.842
$ws.Range("L67").Value = 7082.143
$ws.Range("M67").Value = -2670.842
$ws.Range("N67").Value = -8798.143
$ws.Range("H107").Value = 775.6
$ws.Range("I107").Value = 806
$ws.Range("K107").Value = 806
$ws.Range("M107").Value = 1114
$ws.Range("H112").Value = 33317.656
$ws.Range("J112").Value = 34213.85
$ws.Range("L112").Value = 102641.55
$ws.Range("N112").Value = -104857.55
$ws.Range("H115").Value = 301.875
$ws.Range("I115").Value = 301.875
$ws.Range("K115").Value = 905.625
$ws.Range("M115").Value = 661.375
$ws.Range("H134").Value = 107499.75
$ws.Range("J134").Value = 107499.75
$ws.Range("L134").Value = 107499.75
$ws.Range("N134").Value = -117639.75
$ws.Range("H135").Value = 1372.1923
$ws.Range("I135").Value = 918.6667
$ws.Range("J135").Value = 3277
$ws.Range("K135").Value = 8268.0003
$ws.Range("L135").Value = 29493
$ws.Range("M135").Value = -5733.0003
$ws.Range("N135").Value = -34563
$ws.Range("H137").Value = 1619.921
$ws.Range("I137").Value = 1563.85
$ws.Range("K137").Value = 4691.549999999999
$ws.Range("M137").Value = -2141.549999999999
$ws.Range("H141").Value = 1782.1428
$ws.Range("I141").Value = 1782.1428
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5346.428400000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -166.4284000000007
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6327
$ws.Range("I45").Value = 7123.1177
$ws.Range("J45").Value = 1815.6666
$ws.Range("K45").Value = 7123.1177
$ws.Range("L45").Value = 1815.6666
$ws.Range("M45").Value = -6746.1177
$ws.Range("N45").Value = -2569.6666
$ws.Range("H61").Value = 2863.6287
$ws.Range("I61").Value = 2412.6072
$ws.Range("K61").Value = 2412.6072
$ws.Range("M61").Value = -2200.6072
$ws.Range("H74").Value = 7556.524
$ws.Range("I74").Value = 1767.8
$ws.Range("J74").Value = 12819
$ws.Range("K74").Value = 1767.8
$ws.Range("L74").Value = 12819
$ws.Range("M74").Value = -893.8
$ws.Range("N74").Value = -14567
$ws.Range("H77").Value = 7556.524
$ws.Range("I77").Value = 1767.8
$ws.Range("J77").Value = 12819
$ws.Range("K77").Value = 8839
$ws.Range("L77").Value = 64095
$ws.Range("M77").Value = -4471
$ws.Range("N77").Value = -72831
$ws.Range("H122").Value = 2287.889
$ws.Range("I122").Value = 2170.1428
$ws.Range("K122").Value = 6510.428400000001
$ws.Range("M122").Value = -4060.428400000001
$ws.Range("H132").Value = 3200.4285
$ws.Range("J132").Value = 4292.5557
$ws.Range("L132").Value = 12877.6671
$ws.Range("N132").Value = -17937.6671
$ws.Range("H136").Value = 2863.6287
$ws.Range("I136").Value = 2412.6072
$ws.Range("K136").Value = 7237.821599999999
$ws.Range("M136").Value = -4687.821599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 371727.34
$ws.Range("I22").Value = 686.45
$ws.Range("K22").Value = 686.45
$ws.Range("M22").Value = -513.45
$ws.Range("H130").Value = 85000
$ws.Range("J130").Value = 85000
$ws.Range("L130").Value = 85000
$ws.Range("N130").Value = -95040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 99988
$ws.Range("J127").Value = 99988
$ws.Range("L127").Value = 99988
$ws.Range("N127").Value = -109908
$ws.Range("H135").Value = 69827.14
$ws.Range("J135").Value = 69827.14
$ws.Range("L135").Value = 69827.14
$ws.Range("N135").Value = -79967.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 25001000
$ws.Range("J68").Value = 25001000
$ws.Range("L68").Value = 75003000
$ws.Range("N68").Value = -75004622
$ws.Range("H71").Value = 25001000
$ws.Range("J71").Value = 25001000
$ws.Range("L71").Value = 225009000
$ws.Range("N71").Value = -225017112
$ws.Range("H86").Value = 422.375
$ws.Range("I86").Value = 411.2857
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 1233.8571
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -47.85710000000017
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 422.375
$ws.Range("I89").Value = 411.2857
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 3701.5713
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 2226.4287
$ws.Range("N89").Value = -16356

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H122").Value = 3128.2856
$ws.Range("I122").Value = 2458.0476
$ws.Range("K122").Value = 7374.1428
$ws.Range("M122").Value = -4924.1428
$ws.Range("H132").Value = 3697.7368
$ws.Range("I132").Value = 2903.923
$ws.Range("K132").Value = 8711.769
$ws.Range("M132").Value = -6181.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5444.7334
$ws.Range("I7").Value = 4670.091
$ws.Range("K7").Value = 4670.091
$ws.Range("M7").Value = -4558.091
$ws.Range("H40").Value = 3467.963
$ws.Range("I40").Value = 3151.7646
$ws.Range("J40").Value = 4005.5
$ws.Range("K40").Value = 3151.7646
$ws.Range("L40").Value = 4005.5
$ws.Range("M40").Value = -3015.7646
$ws.Range("N40").Value = -4277.5
$ws.Range("H68").Value = 2840.7693
$ws.Range("I68").Value = 2827.5
$ws.Range("K68").Value = 2827.5
$ws.Range("M68").Value = -2078.5
$ws.Range("H71").Value = 2840.7693
$ws.Range("I71").Value = 2827.5
$ws.Range("K71").Value = 14137.5
$ws.Range("M71").Value = -10393.5
$ws.Range("H122").Value = 5923.1816
$ws.Range("I122").Value = 6159.1665
$ws.Range("J122").Value = 5640
$ws.Range("K122").Value = 18477.4995
$ws.Range("L122").Value = 16920
$ws.Range("M122").Value = -16027.4995
$ws.Range("N122").Value = -21820
$ws.Range("H126").Value = 5444.7334
$ws.Range("I126").Value = 4670.091
$ws.Range("K126").Value = 14010.273
$ws.Range("M126").Value = -11540.273
$ws.Range("H132").Value = 3271.3872
$ws.Range("I132").Value = 3057.08
$ws.Range("J132").Value = 4164.3335
$ws.Range("K132").Value = 9171.24
$ws.Range("L132").Value = 12493.0005
$ws.Range("M132").Value = -6641.24
$ws.Range("N132").Value = -17553.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2120.5715
$ws.Range("I122").Value = 2129.8125
$ws.Range("K122").Value = 6389.4375
$ws.Range("M122").Value = -3939.4375
$ws.Range("H124").Value = 19966.334
$ws.Range("J124").Value = 18449.5
$ws.Range("L124").Value = 18449.5
$ws.Range("N124").Value = -28269.5
$ws.Range("H126").Value = 2533.111
$ws.Range("I126").Value = 2349.75
$ws.Range("K126").Value = 7049.25
$ws.Range("M126").Value = -4579.25
$ws.Range("H133").Value = 47388
$ws.Range("J133").Value = 47388
$ws.Range("L133").Value = 47388
$ws.Range("N133").Value = -57508
